$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw-results row (row 98). Columns A and C contain
# numeric-looking text ("2025-12-23", "251223") that must stay stored as
# text (matching every other row in the sheet), so force a Text number
# format before assigning the values, then restore the default "Normal"
# style so no stray formatting is left behind on the new row.
$row = "98"
$target = $ws.Range("A" + $row + ":E" + $row)
$target.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-12-23"
$ws.Range("B" + $row).Value = "Pick 4"
$ws.Range("C" + $row).Value = "251223"
$ws.Range("D" + $row).Value = "8-7-9-7"
$ws.Range("E" + $row).Value = "2025-12-23T21:42:45.637+04:00"

$target.Style = "Normal"
